$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-style pattern: some Price values look like plain decimal numbers
# (e.g. "298.71"), and Excel's usual text-entry auto-detection would turn
# those into numeric cells (losing the exact string form / trailing zeros
# and introducing float rounding). Forcing NumberFormat to Text ("@")
# before assigning the value keeps it as the literal string, then
# resetting the Style back to "Normal" removes the now-unneeded text
# number-format so the cell's style index matches the rest of the sheet.

# Row 2 - Bitcoin
$ws.Range("D2").Value = "42.822.85"
$ws.Range("E2").Value = "  -0.69%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.289.82"
$ws.Range("E3").Value = "  -1.28%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.04%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "298.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.35%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.77%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +0.45%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.03%  "

# Row 9 - Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.504"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.74%  "

# Row 10 - Avalanche
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.16%  "

# Row 11 - Dogecoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0787"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.57%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +0.86%  "

# Row 13 - Chainlink
$ws.Range("E13").Value = "  -0.36%  "

# Row 14 - Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.12%  "

# Row 15 - Wrapped liquid staked Ether 2.0
$ws.Range("D15").Value = "2.646.98"
$ws.Range("E15").Value = "  -1.22%  "

# Row 16 - Wrapped Ether
$ws.Range("D16").Value = "2.295.40"
$ws.Range("E16").Value = "  -0.94%  "

# Row 17 - Polygon
$ws.Range("E17").Value = "  -2.46%  "

# Row 18 - Wrapped BTC
$ws.Range("D18").Value = "42.747.76"

# Row 19 - Internet Computer (DFINITY)
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.70%  "

# Row 20 - Shiba Inu
$ws.Range("E20").Value = "  -0.56%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  -2.35%  "

# Row 22 - Litecoin
$ws.Range("E22").Value = "  -0.70%  "

# Row 23 - Bitcoin Cash
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.04%  "

# Row 24 - Immutable X
$ws.Range("E24").Value = "  -1.81%  "

# Row 25 - Dai
$ws.Range("E25").Value = "  +0.00%  "

# Row 26 & 27 - PancakeSwap and LEO swap positions
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.44%  "

$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.18%  "

# Row 28 - Ethereum Classic
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.52%  "

# Row 29 - Monero
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "165.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.87%  "

# Row 30 - Toncoin
$ws.Range("E30").Value = "  -0.92%  "

# Row 31 - Cosmos
$ws.Range("E31").Value = "  -1.81%  "

# Row 32 - Injective Protocol
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.96%  "

# Row 33 - First Digital USD
$ws.Range("E33").Value = "  +0.06%  "

# Row 34 - Render Token
$ws.Range("E34").Value = "  -4.71%  "

# Row 35 - Filecoin
$ws.Range("E35").Value = "  -3.81%  "

# Row 36 - Celestia
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.13"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.43%  "

# Row 37 - WEMIX Token
$ws.Range("E37").Value = "  -1.34%  "

# Row 38 - Hedera
$ws.Range("E38").Value = "  -1.68%  "

# Row 39 - Kaspa
$ws.Range("E39").Value = "  -2.03%  "

# Row 40 - ARBITRUM
$ws.Range("E40").Value = "  -4.58%  "

# Row 41 - Lido DAO Token
$ws.Range("E41").Value = "  -1.12%  "

# Row 42 - Stellar
$ws.Range("E42").Value = "  +0.04%  "

# Row 43 - Maker
$ws.Range("D43").Value = "2.013.52"
$ws.Range("E43").Value = "  +0.87%  "

# Row 44 - VeChain
$ws.Range("E44").Value = "  -2.35%  "

# Row 45 - Frax Share
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.06"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.13%  "

# Row 46 - ApeX Protocol
$ws.Range("E46").Value = "  -0.76%  "

# Row 47 - EnergySwap
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.49%  "

# Row 48 - NEAR Protocol
$ws.Range("E48").Value = "  -2.80%  "

# Row 49 - Rocket Pool ETH
$ws.Range("D49").Value = "2.513.69"
$ws.Range("E49").Value = "  -1.25%  "

# Row 50 - MultiversX
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.95"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.45%  "

# Row 51 - HuobiToken -> BitcoinSV
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.62%  "
